$d = $word.ActiveDocument
$d.Content.Find.Execute("Compradora: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Compradora: ", 2)
